$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.599488
$ws.Range("H2").Value = 16.798464
$ws.Range("I2").Value = 0.8205563069596913
$ws.Range("J2").Value = 0.8205563069596913
$ws.Range("M2").Value = 21.106626
$ws.Range("N2").Value = 63.319878
$ws.Range("O2").Value = 0.3106213714361249
$ws.Range("P2").Value = 0.3106213714361249
$ws.Range("Q2").Value = 118.186299007488
$ws.Range("R2").Value = 1063.676691067392
$ws.Range("S2").Value = 0.2548823254083812
$ws.Range("T2").Value = 0.2548823254083812

# Row 3
$ws.Range("G3").Value = 5.599488
$ws.Range("H3").Value = 16.798464
$ws.Range("I3").Value = 0.8205563069596913
$ws.Range("J3").Value = 0.8205563069596913
$ws.Range("O3").Value = 0.03085709917216154
$ws.Range("P3").Value = 0.03085709917216154
$ws.Range("Q3").Value = 11.740616340736
$ws.Range("R3").Value = 105.665547066624
$ws.Range("S3").Value = 0.02531998734019782
$ws.Range("T3").Value = 0.02531998734019782

# Row 4
$ws.Range("G4").Value = 5.599488
$ws.Range("H4").Value = 16.798464
$ws.Range("I4").Value = 0.8205563069596913
$ws.Range("J4").Value = 0.8205563069596913
$ws.Range("O4").Value = 0.6585215293917135
$ws.Range("P4").Value = 0.6585215293917135
$ws.Range("Q4").Value = 250.556560276992
$ws.Range("R4").Value = 2255.009042492928
$ws.Range("S4").Value = 0.5403539942111122
$ws.Range("T4").Value = 0.5403539942111122

# Row 5
$ws.Range("I5").Value = 0.1794436930403087
$ws.Range("J5").Value = 0.1794436930403087
$ws.Range("M5").Value = 21.106626
$ws.Range("N5").Value = 63.319878
$ws.Range("O5").Value = 0.3106213714361249
$ws.Range("P5").Value = 0.3106213714361249
$ws.Range("Q5").Value = 25.845619344818
$ws.Range("R5").Value = 232.610574103362
$ws.Range("S5").Value = 0.05573904602774372
$ws.Range("T5").Value = 0.0557390460277437

# Row 6
$ws.Range("I6").Value = 0.1794436930403087
$ws.Range("J6").Value = 0.1794436930403087
$ws.Range("O6").Value = 0.03085709917216154
$ws.Range("P6").Value = 0.03085709917216154
$ws.Range("S6").Value = 0.005537111831963719
$ws.Range("T6").Value = 0.005537111831963719

# Row 7
$ws.Range("I7").Value = 0.1794436930403087
$ws.Range("J7").Value = 0.1794436930403087
$ws.Range("O7").Value = 0.6585215293917135
$ws.Range("P7").Value = 0.6585215293917135
$ws.Range("S7").Value = 0.1181675351806012
$ws.Range("T7").Value = 0.1181675351806012
